# Gestion de lote - actualizacion de stock (columna H) y precio (columna I)
# para el reporte de productos, y limpieza del resaltado aplicado por error
# a la fila 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# La fila 5 (A5:I5) tenia un estilo "Amount1" (fuente roja) aplicado; se
# revierte al estilo por defecto usado por el resto de filas de datos.
$ws.Range("A5:I5").Font.Color = $ws.Range("A6").Font.Color
$ws.Range("A5:I5").Interior.ColorIndex = 0

# La celda de stock (H5) quedaba vacia; ahora refleja el stock disponible.
$ws.Range("H5").Value = 6

# Actualizacion de stock para los lotes de las filas 6 y 7, y correccion del
# precio de la fila 7.
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 7.6

# La seleccion activa vuelve a la fila de encabezados en vez de la fila 5.
$ws.Range("A4:I4").Select()
